$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-26 12:57:26"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
